$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "09/20/2016 (1hr)" -> "09/20/2016 (2.5hr)", split across two
# bold runs ("09/20/2016 (2.5" / "hr)") as in the target XML.
# ---------------------------------------------------------------------
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$find.Text = "1hr)"
$find.Execute() | Out-Null

# Collapse to just before the "1" and replace it with "2.5".
$rng.Collapse(1)
$rng.MoveEnd(1, 1) | Out-Null
$rng.Text = "2.5"

# Force a run boundary right after "2.5" (before "hr)") by dropping a
# temporary bookmark there and removing it again; adjoining runs with
# identical formatting normally get coalesced back into one <w:r>, but
# while the bookmark sits between them they stay split, and the split
# persists once the bookmark is deleted.
$rng.Collapse(0)
$d.Bookmarks.Add("TempSplitMarker", $rng) | Out-Null
$d.Bookmarks.Item("TempSplitMarker").Delete()

# ---------------------------------------------------------------------
# Change 2: after "Modified code to printing board only after each
# user makes a valid move", add four new bullet paragraphs, then a
# fifth ("NEXT: Serialization and Computer Strategy") split so that the
# _GoBack bookmark (Word's "last edit position" marker) sits between
# "Serializ" and "ation", mirroring where the author's cursor last was.
# ---------------------------------------------------------------------
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$find.Text = "Modified code to printing board only after each user makes a valid move"
$find.Execute() | Out-Null
$rng.Collapse(0)

$texts = @(
    "Implemented the Tournament class and integrated it with the game class properly.",
    "Handled the user choice to quit or continue after each round.",
    "Tested thoroughly to make sure the tournament exits/continues as the user wants",
    "At the point, the game seems to be have all the essential components implemented for it to be a human game.",
    "NEXT: Serialization and Computer Strategy"
)

foreach ($t in $texts) {
    $rng.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last.Range
    $newPara.Collapse(0)
    $newPara.InsertBefore($t)
    $rng = $newPara
    $rng.Collapse(0)
}

# Re-anchor _GoBack between "Serializ" and "ation" in the last
# paragraph ("NEXT: Serializ" is 14 characters).
$lastPara = $d.Paragraphs.Last.Range
$splitPoint = $lastPara.Start + 14
$mid = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $mid) | Out-Null

Write-Output "edit complete"
